# Auto-generated script applying the Universalis price-refresh update
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2349.875
$ws.Range("I80").Value = 700
$ws.Range("J80").Value = 3999.75
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 11999.25
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -13995.25
# Row 83
$ws.Range("H83").Value = 2349.875
$ws.Range("I83").Value = 700
$ws.Range("J83").Value = 3999.75
$ws.Range("K83").Value = 6300
$ws.Range("L83").Value = 35997.75
$ws.Range("M83").Value = -1308
$ws.Range("N83").Value = -45981.75
# Row 116
$ws.Range("H116").Value = 7499
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 7499
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 7499
$ws.Range("N116").Value = -14383
$ws.Range("M116").ClearContents()
# Row 141
$ws.Range("H141").Value = 7949.8335
$ws.Range("I141").Value = 7949.8335
$ws.Range("K141").Value = 23849.5005
$ws.Range("M141").Value = -18669.5005

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 3608.4546
$ws.Range("I88").Value = 1797.6666
$ws.Range("J88").Value = 5781.4
$ws.Range("K88").Value = 1797.6666
$ws.Range("L88").Value = 5781.4
$ws.Range("M88").Value = -1391.6666
$ws.Range("N88").Value = -6593.4
# Row 91
$ws.Range("H91").Value = 3608.4546
$ws.Range("I91").Value = 1797.6666
$ws.Range("J91").Value = 5781.4
$ws.Range("K91").Value = 1797.6666
$ws.Range("L91").Value = 5781.4
$ws.Range("M91").Value = -393.6666
$ws.Range("N91").Value = -8589.4
# Row 92
$ws.Range("H92").Value = 47110
$ws.Range("J92").Value = 47110
$ws.Range("L92").Value = 47110
$ws.Range("N92").Value = -52102
# Row 122
$ws.Range("H122").Value = 1865.5714
$ws.Range("I122").Value = 1553.6666
$ws.Range("K122").Value = 4660.9998
$ws.Range("M122").Value = -2210.9998

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8484
$ws.Range("I86").Value = 6780
$ws.Range("K86").Value = 6780
$ws.Range("M86").Value = -5657
# Row 89
$ws.Range("H89").Value = 8484
$ws.Range("I89").Value = 6780
$ws.Range("K89").Value = 33900
$ws.Range("M89").Value = -28284
# Row 99
$ws.Range("H99").Value = 2712
$ws.Range("I99").Value = 2725
$ws.Range("K99").Value = 2725
$ws.Range("M99").Value = -1227
# Row 107
$ws.Range("H107").Value = 947
$ws.Range("I107").Value = 981.1667
$ws.Range("J107").Value = 844.5
$ws.Range("K107").Value = 981.1667
$ws.Range("L107").Value = 844.5
$ws.Range("M107").Value = 938.8333
$ws.Range("N107").Value = -4684.5
# Row 134
$ws.Range("H134").Value = 2819
$ws.Range("I134").Value = 1430.875
$ws.Range("K134").Value = 4292.625
$ws.Range("M134").Value = -1757.625

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2403.2222
$ws.Range("I31").Value = 1703.625
$ws.Range("K31").Value = 1703.625
$ws.Range("M31").Value = -1408.625
# Row 34
$ws.Range("H34").Value = 2403.2222
$ws.Range("I34").Value = 1703.625
$ws.Range("K34").Value = 1703.625
$ws.Range("M34").Value = -1501.625

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 114.28571
$ws.Range("I7").Value = 125
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 375
$ws.Range("L7").Value = 150
$ws.Range("M7").Value = -263
$ws.Range("N7").Value = -374
# Row 36
$ws.Range("H36").Value = 670.3333
$ws.Range("I36").Value = 20
$ws.Range("J36").Value = 995.5
$ws.Range("K36").Value = 60
$ws.Range("L36").Value = 2986.5
$ws.Range("M36").Value = 109
$ws.Range("N36").Value = -3324.5
# Row 68
$ws.Range("H68").Value = 806.6
$ws.Range("J68").Value = 677.6667
$ws.Range("L68").Value = 2033.0001
$ws.Range("N68").Value = -3655.0001
# Row 71
$ws.Range("H71").Value = 806.6
$ws.Range("J71").Value = 677.6667
$ws.Range("L71").Value = 6099.0003
$ws.Range("N71").Value = -14211.0003
# Row 92
$ws.Range("H92").Value = 787.2857
$ws.Range("I92").Value = 741.8889
$ws.Range("J92").Value = 869
$ws.Range("K92").Value = 2225.6667
$ws.Range("L92").Value = 2607
$ws.Range("M92").Value = -977.6667000000002
$ws.Range("N92").Value = -5103
# Row 97
$ws.Range("H97").Value = 991.5
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
# Row 103
$ws.Range("H103").Value = 1908.3334
$ws.Range("I103").Value = 23
$ws.Range("K103").Value = 69
$ws.Range("M103").Value = 810
# Row 109
$ws.Range("H109").Value = 2268.4
$ws.Range("I109").Value = 335.5
$ws.Range("K109").Value = 1006.5
$ws.Range("M109").Value = 33.5
# Row 122
$ws.Range("H122").Value = 1232
$ws.Range("I122").Value = 1044.4286
$ws.Range("K122").Value = 9399.857399999999
$ws.Range("M122").Value = -6949.857399999999
# Row 131
$ws.Range("H131").Value = 942.8570999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1000
$ws.Range("I80").Value = 1000
$ws.Range("K80").Value = 1000
$ws.Range("M80").Value = -2
# Row 83
$ws.Range("H83").Value = 1000
$ws.Range("I83").Value = 1000
$ws.Range("K83").Value = 5000
$ws.Range("M83").Value = -8
# Row 113
$ws.Range("H113").Value = 849.25
$ws.Range("J113").Value = 899
$ws.Range("L113").Value = 899
$ws.Range("N113").Value = -5239
# Row 126
$ws.Range("H126").Value = 4932.8887
$ws.Range("I126").Value = 2599.6667
$ws.Range("J126").Value = 6099.5
$ws.Range("K126").Value = 7799.000100000001
$ws.Range("L126").Value = 18298.5
$ws.Range("M126").Value = -5329.000100000001
$ws.Range("N126").Value = -23238.5
# Row 132
$ws.Range("H132").Value = 3573.4707
$ws.Range("I132").Value = 3161.5
$ws.Range("K132").Value = 9484.5
$ws.Range("M132").Value = -6954.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2106.6924
$ws.Range("I7").Value = 1829
$ws.Range("J7").Value = 3032.3333
$ws.Range("K7").Value = 1829
$ws.Range("L7").Value = 3032.3333
$ws.Range("M7").Value = -1717
$ws.Range("N7").Value = -3256.3333
# Row 40
$ws.Range("H40").Value = 6099.6875
$ws.Range("I40").Value = 5804.3335
$ws.Range("K40").Value = 5804.3335
$ws.Range("M40").Value = -5668.3335
# Row 126
$ws.Range("H126").Value = 2106.6924
$ws.Range("I126").Value = 1829
$ws.Range("J126").Value = 3032.3333
$ws.Range("K126").Value = 5487
$ws.Range("L126").Value = 9096.999899999999
$ws.Range("M126").Value = -3017
$ws.Range("N126").Value = -14036.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5166.6665
$ws.Range("I62").Value = 5266.6665
$ws.Range("K62").Value = 5266.6665
$ws.Range("M62").Value = -4642.6665
# Row 65
$ws.Range("H65").Value = 5166.6665
$ws.Range("I65").Value = 5266.6665
$ws.Range("K65").Value = 26333.3325
$ws.Range("M65").Value = -23213.3325
# Row 81
$ws.Range("H81").Value = 797.5
$ws.Range("J81").Value = 795
$ws.Range("L81").Value = 1590
$ws.Range("N81").Value = -3712
# Row 84
$ws.Range("H84").Value = 797.5
$ws.Range("J84").Value = 795
$ws.Range("L84").Value = 7950
$ws.Range("N84").Value = -18558
# Row 122
$ws.Range("H122").Value = 706.5625
$ws.Range("I122").Value = 621.2222
$ws.Range("J122").Value = 816.2857
$ws.Range("K122").Value = 1863.6666
$ws.Range("L122").Value = 2448.8571
$ws.Range("M122").Value = 586.3334
$ws.Range("N122").Value = -7348.8571
# Row 126
$ws.Range("H126").Value = 1950
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -10940
# Row 132
$ws.Range("H132").Value = 1804.8235
$ws.Range("I132").Value = 1699.3846
$ws.Range("K132").Value = 5098.1538
$ws.Range("M132").Value = -2568.1538
